$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'94.938.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3
$ws.Range("D3").Value = "'3.614.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'236.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("D6").Value = "'655.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.05%  "

# Row 7
$ws.Range("E7").Value = "  +0.29%  "

# Row 8
$ws.Range("D8").Value = "'0.404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").Value = "'0.991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "

# Row 11
$ws.Range("D11").Value = "'3.613.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.90%  "

# Row 12
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'42.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.40%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.200"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").Value = "'6.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "

# Row 15
$ws.Range("D15").Value = "'4.301.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.03%  "

# Row 16
$ws.Range("D16").Value = "'94.938.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "

# Row 17
$ws.Range("D17").Value = "'0.0000253"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "

# Row 18
$ws.Range("D18").Value = "'3.609.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.94%  "

# Row 19
$ws.Range("D19").Value = "'7.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.49%  "

# Row 20
$ws.Range("D20").Value = "'12.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.59%  "

# Row 21
$ws.Range("D21").Value = "'17.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

# Row 22
$ws.Range("D22").Value = "'3.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.30%  "

# Row 23
$ws.Range("B23").Value = "Stellar"
$ws.Range("C23").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D23").Value = "'0.479"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.04%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'505.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25
$ws.Range("D25").Value = "'0.0000196"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.82%  "

# Row 26
$ws.Range("D26").Value = "'6.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.55%  "

# Row 27
$ws.Range("D27").Value = "'95.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.04%  "

# Row 28
$ws.Range("D28").Value = "'3.807.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "

# Row 29
$ws.Range("D29").Value = "'12.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "

# Row 30
$ws.Range("D30").Value = "'3.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.82%  "

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("D32").Value = "'11.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.22%  "

# Row 33
$ws.Range("D33").Value = "'0.138"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.28%  "

# Row 34
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'32.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.54%  "

# Row 36
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D36").Value = "'0.176"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "

# Row 37
$ws.Range("D37").Value = "'0.557"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "

# Row 38
$ws.Range("D38").Value = "'569.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "

# Row 39
$ws.Range("D39").Value = "'8.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.99%  "

# Row 40
$ws.Range("D40").Value = "'1.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.45%  "

# Row 41
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.916"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.149"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").Value = "'34.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +42.96%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'23.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "

# Row 47
$ws.Range("D47").Value = "'5.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48
$ws.Range("D48").Value = "'2.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.43%  "

# Row 49
$ws.Range("D49").Value = "'0.0412"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.48%  "

# Row 50
$ws.Range("D50").Value = "'3.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.65%  "

# Row 51
$ws.Range("D51").Value = "'53.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
